# TDI Interview Script - final updates
$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# --- Slide 2 body: prepend new sentence before "have" (keeps italic "have" run intact) ---
Replace-Text "I imagine most, if not all of us, have had at least one experience with the flu, probably one we’d be happy to never repeat. And while there exists at least one simple heuristic algorithm for determining if you " "I imagine most, if not all of us, have had at least one experience with the flu, probably one we’d be happy to never repeat. Maybe some of you have seen this diagnosis algorithm floating around on the internet – while this could work to determine if you "

# --- Slide 3 body: drop leading question, rework the middle, swap "national" -> "federal" ---
Replace-Text "Key question for forecasts is what area to forecast for? The Centers for Disease Control and Prevention publishes weekly flu levels for the nation, multi-state regions, and individual states, with the main measure of intensity being the % of doctor visits due to influenza-like illness. If we look at the change in % visits for ILI from late Jan to early Feb nationally, shown in the top left plot, flu activity is going up. But if we look at a regional pattern, shown in the middle, or state pattern, in the bottom right, we see spatial differences suggesting that forecasting at those resolutions would have value for end users. For example, the national government might care about national forecasts, while a pharmaceutical company producing vaccines or antivirals might care about regional or state forecasts to decide on product shipments or marketing." "The Centers for Disease Control and Prevention publishes weekly flu levels for the nation, multi-state regions, and individual states, with the main measure of intensity being the % of doctor visits due to influenza-like illness. At the national level, illustrated in the top left plot, there were more visits for flu in early Feb compared to late Jan. But if we look at a regional pattern, shown in the middle, we can see spatial differences, and zooming in on Region 3 in the Mid-Atlantic shows even more differences between states. This all suggests that forecasting at these varying resolutions would have value for end users. For example, the federal government might care about national forecasts, while a pharmaceutical company producing vaccines or antivirals might care about regional or state forecasts to decide on product shipments or marketing."

# --- Slide 4 body: ILI -> flu virus, reword measure of flu activity, trim trailing phrase ---
Replace-Text "To build my forecasts, I used CDC surveillance data on ILI activity and flu virus types. I used % visits for ILI as my measure of flu activity. These data are collected weekly and available back to 1997 for the US and regions, and back to 2010 for individual states. I also used Google Trends data, available weekly from 2004. The model is based on a dynamic harmonic regression model, where I use Fourier decompositions for the seasonality, ARIMA errors for remaining autocorrelation, and covariates including Google Trends and flu virus types. Forecasts from this model are combined in a weighted ensemble with the results of two naïve baseline models to create the final forecasts. To train the models, I fit weekly forecasts for each CV season to determine the best model structure for each location, which ended up being ~ 14K CV forecasts. From the final forecasts, I estimated the ensemble weights using leave one season out cross validation on the CV seasons." "To build my forecasts, I used CDC surveillance data on flu activity and flu virus types. I used % of doctor visits for flu as my measure of flu activity. These data are collected weekly and available back to 1997 for the US and regions, and back to 2010 for individual states. I also used Google Trends data, available weekly from 2004. Forecasts are based on a dynamic harmonic regression model, where I use Fourier decompositions for the seasonality, ARIMA errors for remaining autocorrelation, and covariates including Google Trends and flu virus types. Forecasts from this model are combined in a weighted ensemble with the results of two naïve baseline models to create the final forecasts. To train the models, I fit weekly forecasts for each CV season to determine the best model structure for each location, which ended up being ~ 14K CV forecasts. From the final forecasts, I estimated the ensemble weights using leave one season out cross validation."

# --- Slide 5 body: two separate edits ---
Replace-Text "I’ve highlighted the performance for the US, Region 3, and two states within Region 3 as examples." "I’ve highlighted the performance for the US, the previously highlighted Region 3, and two states within that as examples."
Replace-Text "levels in PA are forecast to stay fairly steady." "if any of you live in PA you might want to keep washing your hands as levels are predicted to stay steady."

# --- Slide 6 body: two edits (keep the _GoBack bookmark intact, it sits right after this text) ---
Replace-Text "Rather than the static images displayed here, users could click on regions within the map and bring up forecasts for those regions, as well as measures of accuracy. " "Rather than the static mockups displayed here, users could click on regions within the map and bring up forecasts for those regions, "
Replace-Text "By providing a sense of " "providing valuable insights to individuals, health care providers, and public health officials."

# remove the lone trailing space run that followed the _GoBack bookmark
$p12 = $d.Paragraphs.Item(12)
$endPos = $p12.Range.End
$trail = $d.Range($endPos - 2, $endPos - 1)
if ($trail.Text -eq " ") {
    $trail.Text = ""
}

# --- Append trailing paragraphs: blank, "Nation getting smaller", blank, blank ---
$sel = $word.Selection
$sel.EndKey(6)
$sel.TypeParagraph()
$sel.TypeText("Nation getting smaller")
$sel.TypeParagraph()
$sel.TypeParagraph()

# --- Apply 16pt (sz/szCs = 32 half-points) to every paragraph except the new "Nation getting smaller" one ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Nation getting smaller") {
        continue
    }
    $para.Range.Font.Size = 16
    $para.Range.Font.SizeBi = 16
}

Write-Output "done"
